$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 7710.2666
$ws.Range("I9").Value = 8873.462
$ws.Range("J9").Value = 149.5
$ws.Range("K9").Value = 8873.462
$ws.Range("L9").Value = 149.5
$ws.Range("M9").Value = -8704.462
$ws.Range("N9").Value = -487.5
$ws.Range("H33").Value = 675.625
$ws.Range("I33").Value = 715
$ws.Range("K33").Value = 715
$ws.Range("M33").Value = -486
$ws.Range("H51").Value = 5009.364
$ws.Range("I51").Value = 3843
$ws.Range("J51").Value = 7050.5
$ws.Range("K51").Value = 3843
$ws.Range("L51").Value = 7050.5
$ws.Range("M51").Value = -3359
$ws.Range("N51").Value = -8018.5
$ws.Range("H80").Value = 1119.027
$ws.Range("I80").Value = 741.6875
$ws.Range("J80").Value = 1406.5238
$ws.Range("K80").Value = 2225.0625
$ws.Range("L80").Value = 4219.5714
$ws.Range("M80").Value = -1227.0625
$ws.Range("N80").Value = -6215.5714
$ws.Range("H83").Value = 1119.027
$ws.Range("I83").Value = 741.6875
$ws.Range("J83").Value = 1406.5238
$ws.Range("K83").Value = 6675.1875
$ws.Range("L83").Value = 12658.7142
$ws.Range("M83").Value = -1683.1875
$ws.Range("N83").Value = -22642.7142
$ws.Range("H137").Value = 5752.7393
$ws.Range("I137").Value = 5510.6665
$ws.Range("J137").Value = 5908.357
$ws.Range("K137").Value = 16531.9995
$ws.Range("L137").Value = 17725.071
$ws.Range("M137").Value = -13981.9995
$ws.Range("N137").Value = -22825.071
$ws.Range("H138").Value = 8921.52
$ws.Range("I138").Value = 6333.1816
$ws.Range("J138").Value = 10955.214
$ws.Range("K138").Value = 18999.5448
$ws.Range("L138").Value = 32865.642
$ws.Range("M138").Value = -13859.5448
$ws.Range("N138").Value = -43145.642

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6063.1177
$ws.Range("I45").Value = 3156.1667
$ws.Range("J45").Value = 13039.8
$ws.Range("K45").Value = 3156.1667
$ws.Range("L45").Value = 13039.8
$ws.Range("M45").Value = -2779.1667
$ws.Range("N45").Value = -13793.8
$ws.Range("H132").Value = 6986.41
$ws.Range("I132").Value = 2508.6365
$ws.Range("K132").Value = 7525.9095
$ws.Range("M132").Value = -4995.9095

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3225.077
$ws.Range("I99").Value = 2366.875
$ws.Range("K99").Value = 2366.875
$ws.Range("M99").Value = -868.875
$ws.Range("H132").Value = 62499.5
$ws.Range("J132").Value = 62499.5
$ws.Range("L132").Value = 62499.5
$ws.Range("N132").Value = -72619.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5485.2144
$ws.Range("I31").Value = 3732
$ws.Range("J31").Value = 5963.364
$ws.Range("K31").Value = 3732
$ws.Range("L31").Value = 5963.364
$ws.Range("M31").Value = -3437
$ws.Range("N31").Value = -6553.364
$ws.Range("H34").Value = 5485.2144
$ws.Range("I34").Value = 3732
$ws.Range("J34").Value = 5963.364
$ws.Range("K34").Value = 3732
$ws.Range("L34").Value = 5963.364
$ws.Range("M34").Value = -3530
$ws.Range("N34").Value = -6367.364
$ws.Range("H99").Value = 3708.5806
$ws.Range("I99").Value = 2958.7222
$ws.Range("J99").Value = 4746.846
$ws.Range("K99").Value = 2958.7222
$ws.Range("L99").Value = 4746.846
$ws.Range("M99").Value = -1460.7222
$ws.Range("N99").Value = -7742.846
$ws.Range("H107").Value = 1326.862
$ws.Range("I107").Value = 687.95
$ws.Range("J107").Value = 2746.6667
$ws.Range("K107").Value = 687.95
$ws.Range("L107").Value = 2746.6667
$ws.Range("M107").Value = 1232.05
$ws.Range("N107").Value = -6586.6667
$ws.Range("H126").Value = 3708.5806
$ws.Range("I126").Value = 2958.7222
$ws.Range("J126").Value = 4746.846
$ws.Range("K126").Value = 8876.1666
$ws.Range("L126").Value = 14240.538
$ws.Range("M126").Value = -6406.1666
$ws.Range("N126").Value = -19180.538

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 91.625
$ws.Range("I23").Value = 135
$ws.Range("K23").Value = 405
$ws.Range("M23").Value = -170
$ws.Range("H68").Value = 2445.682
$ws.Range("I68").Value = 2252.8948
$ws.Range("J68").Value = 3666.6667
$ws.Range("K68").Value = 6758.6844
$ws.Range("L68").Value = 11000.0001
$ws.Range("M68").Value = -5947.6844
$ws.Range("N68").Value = -12622.0001
$ws.Range("H71").Value = 2445.682
$ws.Range("I71").Value = 2252.8948
$ws.Range("J71").Value = 3666.6667
$ws.Range("K71").Value = 20276.0532
$ws.Range("L71").Value = 33000.0003
$ws.Range("M71").Value = -16220.0532
$ws.Range("N71").Value = -41112.0003
$ws.Range("H87").Value = 16671.334
$ws.Range("I87").Value = 16671.334
$ws.Range("K87").Value = 50014.00199999999
$ws.Range("M87").Value = -48766.00199999999
$ws.Range("H90").Value = 16671.334
$ws.Range("I90").Value = 16671.334
$ws.Range("K90").Value = 150042.006
$ws.Range("M90").Value = -143802.006
$ws.Range("H107").Value = 4424.095
$ws.Range("I107").Value = 1002.8333
$ws.Range("J107").Value = 5792.6
$ws.Range("K107").Value = 3008.4999
$ws.Range("L107").Value = 17377.8
$ws.Range("M107").Value = -1088.4999
$ws.Range("N107").Value = -21217.8
$ws.Range("H109").Value = 125802.375
$ws.Range("I109").Value = 953.1667
$ws.Range("J109").Value = 500350
$ws.Range("K109").Value = 2859.5001
$ws.Range("L109").Value = 1501050
$ws.Range("M109").Value = -1819.5001
$ws.Range("N109").Value = -1503130
$ws.Range("H121").Value = 668635.9
$ws.Range("I121").Value = 1082.5
$ws.Range("J121").Value = 1113671.5
$ws.Range("K121").Value = 3247.5
$ws.Range("L121").Value = 3341014.5
$ws.Range("M121").Value = -1937.5
$ws.Range("N121").Value = -3343634.5
$ws.Range("H131").Value = 3993.6667
$ws.Range("I131").Value = 1426.2858
$ws.Range("J131").Value = 9128.429
$ws.Range("K131").Value = 4278.857400000001
$ws.Range("L131").Value = 27385.287
$ws.Range("M131").Value = 761.1425999999992
$ws.Range("N131").Value = -37465.287

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 60000
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H80").Value = 730894.3
$ws.Range("I80").Value = 491918.28
$ws.Range("J80").Value = 1447822.2
$ws.Range("K80").Value = 491918.28
$ws.Range("L80").Value = 1447822.2
$ws.Range("M80").Value = -490920.28
$ws.Range("N80").Value = -1449818.2
$ws.Range("H83").Value = 730894.3
$ws.Range("I83").Value = 491918.28
$ws.Range("J83").Value = 1447822.2
$ws.Range("K83").Value = 2459591.4
$ws.Range("L83").Value = 7239111
$ws.Range("M83").Value = -2454599.4
$ws.Range("N83").Value = -7249095

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 8280.1
$ws.Range("I61").Value = 6685.857
$ws.Range("J61").Value = 12000
$ws.Range("K61").Value = 6685.857
$ws.Range("L61").Value = 12000
$ws.Range("M61").Value = -6483.857
$ws.Range("N61").Value = -12404
$ws.Range("H113").Value = 8280.1
$ws.Range("I113").Value = 6685.857
$ws.Range("J113").Value = 12000
$ws.Range("K113").Value = 6685.857
$ws.Range("L113").Value = 12000
$ws.Range("M113").Value = -4515.857
$ws.Range("N113").Value = -16340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 20270
$ws.Range("J69").Value = 20270
$ws.Range("L69").Value = 20270
$ws.Range("N69").Value = -21768
$ws.Range("H72").Value = 20270
$ws.Range("J72").Value = 20270
$ws.Range("L72").Value = 60810
$ws.Range("N72").Value = -68298
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H126").Value = 4759.7
$ws.Range("I126").Value = 4678.6
$ws.Range("J126").Value = 4840.8
$ws.Range("K126").Value = 14035.8
$ws.Range("L126").Value = 14522.4
$ws.Range("M126").Value = -11565.8
$ws.Range("N126").Value = -19462.4
$ws.Range("H132").Value = 22544.25
$ws.Range("I132").Value = 4219.357
$ws.Range("K132").Value = 12658.071
$ws.Range("M132").Value = -10128.071
$ws.Range("H133").Value = 64238.332
$ws.Range("J133").Value = 64238.332
$ws.Range("L133").Value = 64238.332
$ws.Range("N133").Value = -74358.33199999999
$ws.Range("H136").Value = 805691
$ws.Range("I136").Value = 1114845.6
$ws.Range("K136").Value = 3344536.8
$ws.Range("M136").Value = -3341986.8
